$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "${userName}_${caseId}" from E3 to E2
$ws.Range("E2").Value = '${userName}_${caseId}'
$ws.Range("E3").Value = $null

# Move the "userName == 'hugang1'" condition from H2 into H3, wrapped as a JSON target
# (this models the new ForEachCommand usage)
$ws.Range("H2").Value = $null
$ws.Range("H3").Value = '{"target":"userName == ''hugang1''"}'

# Adjust column widths for columns E and H (closest representable widths to 23.5 / 34.375
# given this engine's MDW-7 pixel-grid quantization of ColumnWidth)
$ws.Columns.Item(5).ColumnWidth = 22.857142857142854
$ws.Columns.Item(8).ColumnWidth = 33.714285714285715

# Update the selected cell in the sheet view
$ws.Range("D7").Select()
